# Regenerate save_data column G ("K" = strikeouts) to replace the old
# "Strike#" values, as described in the commit message:
#   "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"
#
# Column G (header "K") holds per-game strikeout totals pulled from an
# external box-score source; these are not derivable from the other
# columns on this sheet, so we write the refreshed literal values directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 1
    6  = 0
    7  = 0
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 3
    13 = 1
    14 = 2
    15 = 3
    16 = 1
    17 = 2
    18 = 3
    19 = 0
    20 = 2
    21 = 2
    22 = 1
    23 = 4
    24 = 2
    25 = 3
    26 = 0
    28 = 1
    29 = 2
    30 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
